# Refresh the "cryptos" price table with the latest scrape results.
# Only the cells that actually changed are touched; Coin/Link are only
# set on the rows whose ranking order swapped places.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number => @{ Coin=..; Link=..; Price=..; Volume=.. } (absent keys are left as-is)
$updates = @{
    2 = @{ Price='37.168.15'; Volume='  -0.39%  ' }
    3 = @{ Price='2.005.45'; Volume='  -0.88%  ' }
    4 = @{ Volume='  +0.01%  ' }
    5 = @{ Price='262.70'; Volume='  +6.68%  ' }
    6 = @{ Price='0.610'; Volume='  -1.55%  ' }
    7 = @{ Volume='  +0.21%  ' }
    8 = @{ Price='55.86'; Volume='  -3.43%  ' }
    9 = @{ Price='0.376'; Volume='  -3.04%  ' }
    10 = @{ Price='0.0763'; Volume='  -4.54%  ' }
    11 = @{ Volume='  -3.08%  ' }
    12 = @{ Price='14.21'; Volume='  -4.75%  ' }
    13 = @{ Price='2.290.34'; Volume='  -1.44%  ' }
    14 = @{ Price='22.00'; Volume='  +2.29%  ' }
    15 = @{ Price='0.773'; Volume='  -7.12%  ' }
    16 = @{ Price='5.18'; Volume='  -3.59%  ' }
    17 = @{ Price='2.012.79'; Volume='  -0.79%  ' }
    18 = @{ Price='37.083.12'; Volume='  -0.37%  ' }
    19 = @{ Price='69.76'; Volume='  -0.13%  ' }
    20 = @{ Price='0.0₃0826'; Volume='  -3.01%  ' }
    21 = @{ Price='233.66'; Volume='  +2.66%  ' }
    22 = @{ Price='5.08'; Volume='  -2.25%  ' }
    23 = @{ Volume='  +0.25%  ' }
    24 = @{ Price='2.58'; Volume='  +2.03%  ' }
    25 = @{ Volume='  +0.87%  ' }
    26 = @{ Price='165.34'; Volume='  +1.37%  ' }
    27 = @{ Price='8.86'; Volume='  -3.04%  ' }
    28 = @{ Price='19.37'; Volume='  -2.07%  ' }
    29 = @{ Price='0.128'; Volume='  -6.33%  ' }
    30 = @{ Volume='  -2.85%  ' }
    31 = @{ Price='0.119'; Volume='  -1.56%  ' }
    32 = @{ Price='4.55'; Volume='  -4.20%  ' }
    33 = @{ Price='0.0621'; Volume='  -6.55%  ' }
    34 = @{ Coin='LidoDAOToken'; Link='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; Price='2.42'; Volume='  -1.63%  ' }
    35 = @{ Coin='InternetComputer(DFINITY)'; Link='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; Price='4.36'; Volume='  -4.08%  ' }
    36 = @{ Price='3.51'; Volume='  -0.92%  ' }
    37 = @{ Coin='BinanceUSD'; Link='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; Price='1.00'; Volume='  +0.13%  ' }
    38 = @{ Coin='WEMIXToken'; Link='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; Price='1.80'; Volume='  -0.64%  ' }
    39 = @{ Price='5.38'; Volume='  +0.41%  ' }
    40 = @{ Price='3.06'; Volume='  +1.87%  ' }
    41 = @{ Price='1.18'; Volume='  +1.48%  ' }
    42 = @{ Price='1.442.72'; Volume='  +3.63%  ' }
    43 = @{ Price='0.0912'; Volume='  -5.91%  ' }
    44 = @{ Price='0.0208'; Volume='  -4.09%  ' }
    45 = @{ Price='89.40'; Volume='  -1.24%  ' }
    46 = @{ Price='15.57'; Volume='  -4.77%  ' }
    47 = @{ Price='1.03'; Volume='  -1.31%  ' }
    48 = @{ Price='2.92'; Volume='  +2.02%  ' }
    49 = @{ Price='6.81'; Volume='  -8.40%  ' }
    50 = @{ Price='2.181.86'; Volume='  -1.54%  ' }
    51 = @{ Price='1.93'; Volume='  -7.28%  ' }
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    if ($values.ContainsKey("Coin")) { $ws.Cells.Item($row, 2).Value = $values["Coin"] }
    if ($values.ContainsKey("Link")) { $ws.Cells.Item($row, 3).Value = $values["Link"] }
    if ($values.ContainsKey("Price")) {
        $price = $values["Price"]
        $priceCell = $ws.Cells.Item($row, 4)
        # The Price column holds text like "22.00" or "1.980.32". When the text
        # parses as a plain number (one decimal point, no thousands-style extra
        # dots), Excel silently converts it and drops significant trailing
        # zeros, so force text formatting first in that case.
        if ($price -match '^[0-9]+(\.[0-9]+)?$') {
            $priceCell.NumberFormat = "@"
        }
        $priceCell.Value = $price
    }
    if ($values.ContainsKey("Volume")) { $ws.Cells.Item($row, 5).Value = $values["Volume"] }
}
